# radio bouton facture normale et facture proforma
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update client name (Raison social)
$ws.Range("E2").Value = "Socobis"

# Update invoice line item 18
# A18 holds a text "Bulletin" id (was the text "10"); prefix with an
# apostrophe so Excel keeps it text instead of auto-converting to a number.
$ws.Range("A18").Value = "'20"
$ws.Range("B18").Value = "Frego"
$ws.Range("D18").Value = 150000
$ws.Range("E18").Value = 240000
$ws.Range("G18").Value = 390000
